{"js": "// Positional old->new text replacements, in document order (date paragraph\n// followed by every table-cell paragraph). Duplicate \"old\" values exist\n// (e.g. \"60-42=18\" and \"72-2=70\" each occur twice with different\n// replacements), so cells are matched by index, not by text search.\nconst REPLACEMENTS = [\n  [\"2025-09-17 Wednesday\", \"2025-09-18 Thursday\"],\n  [\"53+2=55\", \"45-40=5\"],\n  [\"49+42=91\", \"80+18=98\"],\n  [\"86-25=61\", \"70+19=89\"],\n  [\"40+41=81\", \"7-4=3\"],\n  [\"82-49=33\", \"3+15=18\"],\n  [\"13+74=87\", \"69-50=19\"],\n  [\"29-27=2\", \"41-10=31\"],\n  [\"35+49=84\", \"83-62=21\"],\n  [\"50-23=27\", \"64-30=34\"],\n  [\"44-8=36\", \"99-23=76\"],\n  [\"9+57=66\", \"29+27=56\"],\n  [\"72-42=30\", \"88-70=18\"],\n  [\"25+7=32\", \"84-26=58\"],\n  [\"69+11=80\", \"15+35=50\"],\n  [\"23-15=8\", \"14-10=4\"],\n  [\"15+39=54\", \"19+61=80\"],\n  [\"3+5=8\", \"17+43=60\"],\n  [\"60-42=18\", \"24+13=37\"],\n  [\"85-14=71\", \"73-25=48\"],\n  [\"91-29=62\", \"65-44=21\"],\n  [\"39+23=62\", \"45-34=11\"],\n  [\"24+38=62\", \"43+32=75\"],\n  [\"2+53=55\", \"77-38=39\"],\n  [\"52+47=99\", \"29+10=39\"],\n  [\"18+59=77\", \"84-25=59\"],\n  [\"92-27=65\", \"12+84=96\"],\n  [\"12-7=5\", \"42+1=43\"],\n  [\"45-29=16\", \"29+21=50\"],\n  [\"71-39=32\", \"29+40=69\"],\n  [\"58-46=12\", \"78-13=65\"],\n  [\"10+50=60\", \"53+33=86\"],\n  [\"32+62=94\", \"64+35=99\"],\n  [\"54-17=37\", \"60+31=91\"],\n  [\"24-17=7\", \"15+66=81\"],\n  [\"46-17=29\", \"0+41=41\"],\n  [\"93-56=37\", \"12+52=64\"],\n  [\"37-4=33\", \"13+68=81\"],\n  [\"82-9=73\", \"72-58=14\"],\n  [\"65-4=61\", \"50-12=38\"],\n  [\"37-34=3\", \"59-45=14\"],\n  [\"5+21=26\", \"35-31=4\"],\n  [\"73+21=94\", \"19+58=77\"],\n  [\"47+46=93\", \"58-56=2\"],\n  [\"45+46=91\", \"29+53=82\"],\n  [\"34-25=9\", \"75+10=85\"],\n  [\"67+4=71\", \"96-82=14\"],\n  [\"36+33=69\", \"28+28=56\"],\n  [\"76-10=66\", \"89-9=80\"],\n  [\"14+50=64\", \"66-64=2\"],\n  [\"17+51=68\", \"15+52=67\"],\n  [\"98-62=36\", \"38-28=10\"],\n  [\"74-30=44\", \"69-49=20\"],\n  [\"32+7=39\", \"83-81=2\"],\n  [\"32+46=78\", \"8+23=31\"],\n  [\"66+13=79\", \"89-78=11\"],\n  [\"34-14=20\", \"37+30=67\"],\n  [\"51-33=18\", \"56-9=47\"],\n  [\"15+22=37\", \"14+79=93\"],\n  [\"42+2=44\", \"20-0=20\"],\n  [\"44+8=52\", \"40-20=20\"],\n  [\"39+28=67\", \"90-74=16\"],\n  [\"5+1=6\", \"75-57=18\"],\n  [\"20+47=67\", \"18+33=51\"],\n  [\"31-15=16\", \"71-61=10\"],\n  [\"82+1=83\", \"94-19=75\"],\n  [\"34+38=72\", \"81-51=30\"],\n  [\"72-2=70\", \"21-5=16\"],\n  [\"22+67=89\", \"14+13=27\"],\n  [\"74-26=48\", \"31-24=7\"],\n  [\"92-45=47\", \"51+45=96\"],\n  [\"95-19=76\", \"33+15=48\"],\n  [\"70-23=47\", \"70+19=89\"],\n  [\"92-0=92\", \"87-64=23\"],\n  [\"70-12=58\", \"9+6=15\"],\n  [\"47-13=34\", \"79-59=20\"],\n  [\"61+32=93\", \"79-30=49\"],\n  [\"64-43=21\", \"60-12=48\"],\n  [\"59-18=41\", \"50-49=1\"],\n  [\"60+37=97\", \"54+0=54\"],\n  [\"59-44=15\", \"2+52=54\"],\n  [\"85+1=86\", \"46+10=56\"],\n  [\"42+33=75\", \"17+55=72\"],\n  [\"54+13=67\", \"84+12=96\"],\n  [\"22+31=53\", \"32+26=58\"],\n  [\"60-42=18\", \"6+23=29\"],\n  [\"11+17=28\", \"29+38=67\"],\n  [\"3+93=96\", \"48+51=99\"],\n  [\"1+14=15\", \"78-11=67\"],\n  [\"21+59=80\", \"54-10=44\"],\n  [\"69-42=27\", \"62+8=70\"],\n  [\"29+3=32\", \"0+32=32\"],\n  [\"72-2=70\", \"63-34=29\"],\n  [\"8+22=30\", \"19+42=61\"],\n  [\"15+14=29\", \"83-51=32\"],\n  [\"20-11=9\", \"5+71=76\"],\n  [\"89-24=65\", \"10+15=25\"],\n  [\"92-63=29\", \"63-55=8\"],\n  [\"28-23=5\", \"57+31=88\"],\n  [\"23+71=94\", \"13-9=4\"],\n  [\"29+37=66\", \"61-26=35\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nif (items.length !== REPLACEMENTS.length) {\n  throw new Error(\n    `Expected ${REPLACEMENTS.length} paragraphs, found ${items.length}`\n  );\n}\n\nfor (let i = 0; i < items.length; i++) {\n  const [oldText, newText] = REPLACEMENTS[i];\n  const para = items[i];\n  if (para.text !== oldText) {\n    throw new Error(\n      `Paragraph ${i}: expected \"${oldText}\", found \"${para.text}\"`\n    );\n  }\n  if (oldText !== newText) {\n    para.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Positional old->new text replacements, in document order (date paragraph\n# followed by every table-cell paragraph). Duplicate \"old\" values exist\n# (e.g. \"60-42=18\" and \"72-2=70\" each occur twice with different\n# replacements), so cells are matched by index, not by text search.\n$replacements = @(\n    @('2025-09-17 Wednesday', '2025-09-18 Thursday'),\n    @('53+2=55', '45-40=5'),\n    @('49+42=91', '80+18=98'),\n    @('86-25=61', '70+19=89'),\n    @('40+41=81', '7-4=3'),\n    @('82-49=33', '3+15=18'),\n    @('13+74=87', '69-50=19'),\n    @('29-27=2', '41-10=31'),\n    @('35+49=84', '83-62=21'),\n    @('50-23=27', '64-30=34'),\n    @('44-8=36', '99-23=76'),\n    @('9+57=66', '29+27=56'),\n    @('72-42=30', '88-70=18'),\n    @('25+7=32', '84-26=58'),\n    @('69+11=80', '15+35=50'),\n    @('23-15=8', '14-10=4'),\n    @('15+39=54', '19+61=80'),\n    @('3+5=8', '17+43=60'),\n    @('60-42=18', '24+13=37'),\n    @('85-14=71', '73-25=48'),\n    @('91-29=62', '65-44=21'),\n    @('39+23=62', '45-34=11'),\n    @('24+38=62', '43+32=75'),\n    @('2+53=55', '77-38=39'),\n    @('52+47=99', '29+10=39'),\n    @('18+59=77', '84-25=59'),\n    @('92-27=65', '12+84=96'),\n    @('12-7=5', '42+1=43'),\n    @('45-29=16', '29+21=50'),\n    @('71-39=32', '29+40=69'),\n    @('58-46=12', '78-13=65'),\n    @('10+50=60', '53+33=86'),\n    @('32+62=94', '64+35=99'),\n    @('54-17=37', '60+31=91'),\n    @('24-17=7', '15+66=81'),\n    @('46-17=29', '0+41=41'),\n    @('93-56=37', '12+52=64'),\n    @('37-4=33', '13+68=81'),\n    @('82-9=73', '72-58=14'),\n    @('65-4=61', '50-12=38'),\n    @('37-34=3', '59-45=14'),\n    @('5+21=26', '35-31=4'),\n    @('73+21=94', '19+58=77'),\n    @('47+46=93', '58-56=2'),\n    @('45+46=91', '29+53=82'),\n    @('34-25=9', '75+10=85'),\n    @('67+4=71', '96-82=14'),\n    @('36+33=69', '28+28=56'),\n    @('76-10=66', '89-9=80'),\n    @('14+50=64', '66-64=2'),\n    @('17+51=68', '15+52=67'),\n    @('98-62=36', '38-28=10'),\n    @('74-30=44', '69-49=20'),\n    @('32+7=39', '83-81=2'),\n    @('32+46=78', '8+23=31'),\n    @('66+13=79', '89-78=11'),\n    @('34-14=20', '37+30=67'),\n    @('51-33=18', '56-9=47'),\n    @('15+22=37', '14+79=93'),\n    @('42+2=44', '20-0=20'),\n    @('44+8=52', '40-20=20'),\n    @('39+28=67', '90-74=16'),\n    @('5+1=6', '75-57=18'),\n    @('20+47=67', '18+33=51'),\n    @('31-15=16', '71-61=10'),\n    @('82+1=83', '94-19=75'),\n    @('34+38=72', '81-51=30'),\n    @('72-2=70', '21-5=16'),\n    @('22+67=89', '14+13=27'),\n    @('74-26=48', '31-24=7'),\n    @('92-45=47', '51+45=96'),\n    @('95-19=76', '33+15=48'),\n    @('70-23=47', '70+19=89'),\n    @('92-0=92', '87-64=23'),\n    @('70-12=58', '9+6=15'),\n    @('47-13=34', '79-59=20'),\n    @('61+32=93', '79-30=49'),\n    @('64-43=21', '60-12=48'),\n    @('59-18=41', '50-49=1'),\n    @('60+37=97', '54+0=54'),\n    @('59-44=15', '2+52=54'),\n    @('85+1=86', '46+10=56'),\n    @('42+33=75', '17+55=72'),\n    @('54+13=67', '84+12=96'),\n    @('22+31=53', '32+26=58'),\n    @('60-42=18', '6+23=29'),\n    @('11+17=28', '29+38=67'),\n    @('3+93=96', '48+51=99'),\n    @('1+14=15', '78-11=67'),\n    @('21+59=80', '54-10=44'),\n    @('69-42=27', '62+8=70'),\n    @('29+3=32', '0+32=32'),\n    @('72-2=70', '63-34=29'),\n    @('8+22=30', '19+42=61'),\n    @('15+14=29', '83-51=32'),\n    @('20-11=9', '5+71=76'),\n    @('89-24=65', '10+15=25'),\n    @('92-63=29', '63-55=8'),\n    @('28-23=5', '57+31=88'),\n    @('23+71=94', '13-9=4'),\n    @('29+37=66', '61-26=35'),\n)\n\n$d = $word.ActiveDocument\n$total = $d.Paragraphs.Count\n\n# Word's COM paragraph model inserts a synthetic end-of-row marker\n# paragraph (text = CR + cell-mark, i.e. length 2) after every table row's\n# last cell; skip those so only the real content paragraphs remain, in\n# document order, which lines up 1:1 with $replacements.\n$idx = 0\nfor ($i = 1; $i -le $total; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $r = $p.Range\n    $full = $r.Text\n    $trimmed = $full.TrimEnd([char]13, [char]7)\n    if ($trimmed -eq \"\") {\n        continue\n    }\n    if ($idx -ge $replacements.Count) {\n        throw \"More content paragraphs than replacements\"\n    }\n    $oldText = $replacements[$idx][0]\n    $newText = $replacements[$idx][1]\n    if ($trimmed -ne $oldText) {\n        throw \"Paragraph $i`: expected [$oldText], found [$trimmed]\"\n    }\n    if ($oldText -ne $newText) {\n        $r.Text = $newText\n    }\n    $idx++\n}\n\nif ($idx -ne $replacements.Count) {\n    throw \"Expected $($replacements.Count) content paragraphs, found $idx\"\n}\n"}
